# Generate Report for Handoff
# Adds two new tracked files (b4842130-...md and its dependent .png assets)
# to the localization-status workbook: Overview sheet plus the zh-cn / de-de
# detail sheets. The existing "e2bfda48-...md" source file is replaced by
# the new "30564267-...png" / "b9643a17-...png" / "b4842130-...md" trio, and
# the 4th data row (".localization-config") shifts down to row 5.

$wb = $excel.ActiveWorkbook

$linkColor = 15570276   # BGR for RGB(100,149,237) == FF6495ED, matches the
                         # workbook's existing HyperLink cell style.

function Style-AsLink($range) {
    # Applied AFTER Hyperlinks.Add (which stamps its own theme-colored
    # font) so the workbook's custom blue survives in the saved style.
    $range.Font.Underline = 2     # xlUnderlineStyleSingle
    $range.Font.Color = $linkColor
}

function Style-AsDate($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value = "30564267-c082-48d5-8328-00ea5882e309.png"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("A3").Value = "b4842130-f142-49f0-aa76-bcb44163420c.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A4").Value = "b9643a17-5bbf-46f0-9d1d-d4c623a5cb58.png"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

$ws1.Range("A5").Value = ".localization-config"
$ws1.Range("B5").Value = "Not to be localized"
$ws1.Range("C5").Value = "Not to be localized"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6214f9401dab107754b875c134180abb5d387cb0/e2e/30564267-c082-48d5-8328-00ea5882e309.png", "", "", "30564267-c082-48d5-8328-00ea5882e309.png")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6214f9401dab107754b875c134180abb5d387cb0/e2e/b4842130-f142-49f0-aa76-bcb44163420c.md", "", "", "b4842130-f142-49f0-aa76-bcb44163420c.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6214f9401dab107754b875c134180abb5d387cb0/e2e/b9643a17-5bbf-46f0-9d1d-d4c623a5cb58.png", "", "", "b9643a17-5bbf-46f0-9d1d-d4c623a5cb58.png")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/6214f9401dab107754b875c134180abb5d387cb0/.localization-config", "", "", ".localization-config")

Style-AsLink($ws1.Range("A2:A5"))

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("A2").Value = "30564267-c082-48d5-8328-00ea5882e309.png"
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "5a45090c2d2e60e60d153feb5319694dcab6e23d.png"
$ws2.Range("D2").Value = "2016-03-09 01:30:40"
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "IsDependency"
$ws2.Range("I2").Value = "e2e\b4842130-f142-49f0-aa76-bcb44163420c.md"

$ws2.Range("A3").Value = "b4842130-f142-49f0-aa76-bcb44163420c.md"
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = "b4842130-f142-49f0-aa76-bcb44163420c.6b3c89a8b1bf6bd2aa8fc3435211161fa69b3a0d.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-03-09 01:30:40"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = "b9643a17-5bbf-46f0-9d1d-d4c623a5cb58.png"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "1148297276b8751058fc068d29129ebf9922ce1d.png"
$ws2.Range("D4").Value = "2016-03-09 01:30:40"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "IsDependency"
$ws2.Range("I4").Value = "e2e\b4842130-f142-49f0-aa76-bcb44163420c.md"

$ws2.Range("A5").Value = ".localization-config"
$ws2.Range("B5").Value = "Not to be localized"
$ws2.Range("D5").Value = "0001-01-01 00:00:00"
$ws2.Range("G5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6214f9401dab107754b875c134180abb5d387cb0/e2e/30564267-c082-48d5-8328-00ea5882e309.png", "", "", "30564267-c082-48d5-8328-00ea5882e309.png")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/93082246fc73893bdcc4eb79a2dfdbdc14d3343c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5a45090c2d2e60e60d153feb5319694dcab6e23d.png", "", "", "5a45090c2d2e60e60d153feb5319694dcab6e23d.png")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6214f9401dab107754b875c134180abb5d387cb0/e2e/b4842130-f142-49f0-aa76-bcb44163420c.md", "", "", "b4842130-f142-49f0-aa76-bcb44163420c.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/93082246fc73893bdcc4eb79a2dfdbdc14d3343c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b4842130-f142-49f0-aa76-bcb44163420c.6b3c89a8b1bf6bd2aa8fc3435211161fa69b3a0d.zh-cn.xlf", "", "", "b4842130-f142-49f0-aa76-bcb44163420c.6b3c89a8b1bf6bd2aa8fc3435211161fa69b3a0d.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6214f9401dab107754b875c134180abb5d387cb0/e2e/b9643a17-5bbf-46f0-9d1d-d4c623a5cb58.png", "", "", "b9643a17-5bbf-46f0-9d1d-d4c623a5cb58.png")
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/93082246fc73893bdcc4eb79a2dfdbdc14d3343c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/1148297276b8751058fc068d29129ebf9922ce1d.png", "", "", "1148297276b8751058fc068d29129ebf9922ce1d.png")
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/6214f9401dab107754b875c134180abb5d387cb0/.localization-config", "", "", ".localization-config")

Style-AsLink($ws2.Range("A2:A5"))
Style-AsLink($ws2.Range("C2:C4"))
Style-AsDate($ws2.Range("D2:D5"))

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("A2").Value = "30564267-c082-48d5-8328-00ea5882e309.png"
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "5a45090c2d2e60e60d153feb5319694dcab6e23d.png"
$ws3.Range("D2").Value = "2016-03-09 01:30:50"
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "IsDependency"
$ws3.Range("I2").Value = "e2e\b4842130-f142-49f0-aa76-bcb44163420c.md"

$ws3.Range("A3").Value = "b4842130-f142-49f0-aa76-bcb44163420c.md"
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = "b4842130-f142-49f0-aa76-bcb44163420c.6b3c89a8b1bf6bd2aa8fc3435211161fa69b3a0d.de-de.xlf"
$ws3.Range("D3").Value = "2016-03-09 01:30:50"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = "b9643a17-5bbf-46f0-9d1d-d4c623a5cb58.png"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "1148297276b8751058fc068d29129ebf9922ce1d.png"
$ws3.Range("D4").Value = "2016-03-09 01:30:50"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "IsDependency"
$ws3.Range("I4").Value = "e2e\b4842130-f142-49f0-aa76-bcb44163420c.md"

$ws3.Range("A5").Value = ".localization-config"
$ws3.Range("B5").Value = "Not to be localized"
$ws3.Range("D5").Value = "0001-01-01 00:00:00"
$ws3.Range("G5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6214f9401dab107754b875c134180abb5d387cb0/e2e/30564267-c082-48d5-8328-00ea5882e309.png", "", "", "30564267-c082-48d5-8328-00ea5882e309.png")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cbcb0df7de922f179544b4aba8c09843f1efa982/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5a45090c2d2e60e60d153feb5319694dcab6e23d.png", "", "", "5a45090c2d2e60e60d153feb5319694dcab6e23d.png")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6214f9401dab107754b875c134180abb5d387cb0/e2e/b4842130-f142-49f0-aa76-bcb44163420c.md", "", "", "b4842130-f142-49f0-aa76-bcb44163420c.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cbcb0df7de922f179544b4aba8c09843f1efa982/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b4842130-f142-49f0-aa76-bcb44163420c.6b3c89a8b1bf6bd2aa8fc3435211161fa69b3a0d.de-de.xlf", "", "", "b4842130-f142-49f0-aa76-bcb44163420c.6b3c89a8b1bf6bd2aa8fc3435211161fa69b3a0d.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6214f9401dab107754b875c134180abb5d387cb0/e2e/b9643a17-5bbf-46f0-9d1d-d4c623a5cb58.png", "", "", "b9643a17-5bbf-46f0-9d1d-d4c623a5cb58.png")
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cbcb0df7de922f179544b4aba8c09843f1efa982/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/1148297276b8751058fc068d29129ebf9922ce1d.png", "", "", "1148297276b8751058fc068d29129ebf9922ce1d.png")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/6214f9401dab107754b875c134180abb5d387cb0/.localization-config", "", "", ".localization-config")

Style-AsLink($ws3.Range("A2:A5"))
Style-AsLink($ws3.Range("C2:C4"))
Style-AsDate($ws3.Range("D2:D5"))
